{"js": "// Office.js (Word JavaScript API) edit script\n//\n// Target change (see commit message \"Rechtschreibfehler korrigiert\n// Niedergrie\u00dfe -> Niedergriese\"):\n//  1) Remove the \"_GoBack\" bookmark that used to sit right after the title\n//     text \" METLAB\" (end of the first paragraph) -- a leftover cursor\n//     marker from a previous edit.\n//  2) Fix the spelling \"Niedergrie\u00dfe\" -> \"Niedergriese\" inside the role\n//     table (\"Implementierung\" row), and leave a \"_GoBack\" bookmark\n//     positioned right before the final \"e\" of the corrected name (i.e.\n//     Word's live-edit cursor marker lands where the correction was typed).\n\n// --- Step 1: drop the old _GoBack bookmark, wherever it currently is -------\nconst oldMark = context.document.getBookmarkRangeOrNullObject(\"_GoBack\");\noldMark.load(\"isNullObject\");\nawait context.sync();\nif (!oldMark.isNullObject) {\n  context.document.deleteBookmark(\"_GoBack\");\n  await context.sync();\n}\n\n// --- Step 2: correct the misspelled surname ---------------------------------\nconst misspelled = context.document.body.search(\"Niedergrie\u00dfe\", { matchCase: true });\nmisspelled.load(\"items/text\");\nawait context.sync();\n\nif (misspelled.items.length > 0) {\n  // \"\u00df\" + trailing \"e\" becomes \"s\" + \"e\" -> \"Niedergriese\".\n  misspelled.items[0].insertText(\"Niedergriese\", \"Replace\");\n  await context.sync();\n}\n\n// --- Step 3: re-insert \"_GoBack\" right before the final \"e\" of the fix -----\n// This mirrors where Word leaves its \"last edit\" bookmark after a live\n// correction: right after \"...Niedergries\" and before the closing \"e\".\nconst fixedPrefix = context.document.body.search(\"Niedergries\", { matchCase: true });\nfixedPrefix.load(\"items\");\nawait context.sync();\n\nif (fixedPrefix.items.length > 0) {\n  const caret = fixedPrefix.items[0].getRange(\"After\");\n  caret.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# Word COM interop edit script\n#\n# Target change (see commit message \"Rechtschreibfehler korrigiert\n# Niedergrie\u00dfe -> Niedergriese\"):\n#  1) Remove the \"_GoBack\" bookmark that used to sit right after the title\n#     text \" METLAB\" (end of the first paragraph) -- a leftover cursor\n#     marker from a previous edit.\n#  2) Fix the spelling \"Niedergrie\u00dfe\" -> \"Niedergriese\" inside the role\n#     table (\"Implementierung\" row), and leave a \"_GoBack\" bookmark\n#     positioned right before the final \"e\" of the corrected name (i.e.\n#     Word's live-edit cursor marker lands where the correction was typed).\n\n$d = $word.ActiveDocument\n\n# --- Step 1: drop the old _GoBack bookmark, wherever it currently is -------\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# --- Step 2: correct the misspelled surname ---------------------------------\n$findRange = $d.Content\n$findRange.Find.ClearFormatting()\n$findRange.Find.Execute(\"Niedergrie\u00dfe\", $false, $false, $false, $false, $false, $true, 1, $false, \"Niedergriese\", 2)\n\n# --- Step 3: re-insert \"_GoBack\" right before the final \"e\" of the fix -----\n# This mirrors where Word leaves its \"last edit\" bookmark after a live\n# correction: right after \"...Niedergries\" and before the closing \"e\".\n$caretRange = $d.Content\nif ($caretRange.Find.Execute(\"Niedergries\")) {\n    $caretRange.Collapse(0)\n    $d.Bookmarks.Add(\"_GoBack\", $caretRange)\n}\n"}
